# Auto-generated Excel COM-interop edit script
# Applies cell-value updates (refreshed market-price data) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets of the workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 725
$ws.Range("J17").Value = 666.6667
$ws.Range("L17").Value = 2000.0001
$ws.Range("N17").Value = -2336.0001

$ws.Range("H41").Value = 371
$ws.Range("I41").Value = 345.2
$ws.Range("K41").Value = 345.2
$ws.Range("M41").Value = 94.80000000000001

$ws.Range("H55").Value = 265.625
$ws.Range("I55").Value = 178
$ws.Range("K55").Value = 178
$ws.Range("M55").Value = 36

$ws.Range("H70").Value = 1809.7778
$ws.Range("J70").Value = 2042.1428
$ws.Range("L70").Value = 6126.428400000001
$ws.Range("N70").Value = -6666.428400000001

$ws.Range("H73").Value = 1809.7778
$ws.Range("J73").Value = 2042.1428
$ws.Range("L73").Value = 6126.428400000001
$ws.Range("N73").Value = -7998.428400000001

$ws.Range("H80").Value = 4218.75
$ws.Range("I80").Value = 3375
$ws.Range("J80").Value = 5062.5
$ws.Range("K80").Value = 10125
$ws.Range("L80").Value = 15187.5
$ws.Range("M80").Value = -9127
$ws.Range("N80").Value = -17183.5

$ws.Range("H83").Value = 4218.75
$ws.Range("I83").Value = 3375
$ws.Range("J83").Value = 5062.5
$ws.Range("K83").Value = 30375
$ws.Range("L83").Value = 45562.5
$ws.Range("M83").Value = -25383
$ws.Range("N83").Value = -55546.5

$ws.Range("H98").Value = 917.875
$ws.Range("J98").Value = 999
$ws.Range("L98").Value = 999
$ws.Range("N98").Value = -3995

$ws.Range("H122").Value = 917.875
$ws.Range("J122").Value = 999
$ws.Range("L122").Value = 2997
$ws.Range("N122").Value = -7897

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2971.4443
$ws.Range("I88").Value = 2031
$ws.Range("J88").Value = 3441.6667
$ws.Range("K88").Value = 2031
$ws.Range("L88").Value = 3441.6667
$ws.Range("M88").Value = -1625
$ws.Range("N88").Value = -4253.6667

$ws.Range("H91").Value = 2971.4443
$ws.Range("I91").Value = 2031
$ws.Range("J91").Value = 3441.6667
$ws.Range("K91").Value = 2031
$ws.Range("L91").Value = 3441.6667
$ws.Range("M91").Value = -627
$ws.Range("N91").Value = -6249.6667

$ws.Range("H122").Value = 799.6667
$ws.Range("I122").Value = 799.6667
$ws.Range("K122").Value = 2399.0001
$ws.Range("M122").Value = 50.9998999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4665.4165
$ws.Range("J94").Value = 4784.5713
$ws.Range("L94").Value = 4784.5713
$ws.Range("N94").Value = -5686.5713

$ws.Range("H99").Value = 3575.7144
$ws.Range("I99").Value = 3575.7144
$ws.Range("K99").Value = 3575.7144
$ws.Range("M99").Value = -2077.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 471.76923
$ws.Range("I22").Value = 478.3
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 478.3
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = -128.3
$ws.Range("N22").Value = -1150

$ws.Range("H105").Value = 21338.8
$ws.Range("I105").Value = 26173.5
$ws.Range("K105").Value = 26173.5
$ws.Range("M105").Value = -24426.5

$ws.Range("H107").Value = 400.44446
$ws.Range("I107").Value = 374
$ws.Range("J107").Value = 433.5
$ws.Range("K107").Value = 374
$ws.Range("L107").Value = 433.5
$ws.Range("M107").Value = 1546
$ws.Range("N107").Value = -4273.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2203.2727
$ws.Range("I129").Value = 869.3333
$ws.Range("J129").Value = 2703.5
$ws.Range("K129").Value = 2607.9999
$ws.Range("L129").Value = 8110.5
$ws.Range("M129").Value = 2392.0001
$ws.Range("N129").Value = -18110.5

$ws.Range("H140").Value = 752.4
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4289.3335
$ws.Range("J80").Value = 5953
$ws.Range("L80").Value = 5953
$ws.Range("N80").Value = -7949

$ws.Range("H83").Value = 4289.3335
$ws.Range("J83").Value = 5953
$ws.Range("L83").Value = 29765
$ws.Range("N83").Value = -39749

$ws.Range("H122").Value = 1775
$ws.Range("I122").Value = 1775
$ws.Range("K122").Value = 5325
$ws.Range("M122").Value = -2875

$ws.Range("H126").Value = 5300
$ws.Range("I126").Value = 4950
$ws.Range("K126").Value = 14850
$ws.Range("M126").Value = -12380

$ws.Range("H132").Value = 4464.2856
$ws.Range("I132").Value = 5300
$ws.Range("J132").Value = 2375
$ws.Range("K132").Value = 15900
$ws.Range("L132").Value = 7125
$ws.Range("M132").Value = -13370
$ws.Range("N132").Value = -12185

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4999.5
$ws.Range("J7").Value = 4999.5
$ws.Range("L7").Value = 4999.5
$ws.Range("N7").Value = -5223.5

$ws.Range("H22").Value = 2827.4614
$ws.Range("I22").Value = 498
$ws.Range("J22").Value = 3526.3
$ws.Range("K22").Value = 498
$ws.Range("L22").Value = 3526.3
$ws.Range("M22").Value = -203
$ws.Range("N22").Value = -4116.3

$ws.Range("H27").Value = 2827.4614
$ws.Range("I27").Value = 498
$ws.Range("J27").Value = 3526.3
$ws.Range("K27").Value = 498
$ws.Range("L27").Value = 3526.3
$ws.Range("M27").Value = -391
$ws.Range("N27").Value = -3740.3

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H68").Value = 1750
$ws.Range("I68").Value = 1000
$ws.Range("K68").Value = 1000
$ws.Range("M68").Value = -251

$ws.Range("H71").Value = 1750
$ws.Range("I71").Value = 1000
$ws.Range("K71").Value = 5000
$ws.Range("M71").Value = -1256

$ws.Range("H82").Value = 577.5454999999999
$ws.Range("I82").Value = 597.5
$ws.Range("J82").Value = 553.6
$ws.Range("K82").Value = 597.5
$ws.Range("L82").Value = 553.6
$ws.Range("M82").Value = -236.5
$ws.Range("N82").Value = -1275.6

$ws.Range("H85").Value = 577.5454999999999
$ws.Range("I85").Value = 597.5
$ws.Range("J85").Value = 553.6
$ws.Range("K85").Value = 597.5
$ws.Range("L85").Value = 553.6
$ws.Range("M85").Value = 650.5
$ws.Range("N85").Value = -3049.6

$ws.Range("H122").Value = 10000
$ws.Range("I122").Value = 10000
$ws.Range("K122").Value = 30000
$ws.Range("M122").Value = -27550

$ws.Range("H126").Value = 4999.5
$ws.Range("J126").Value = 4999.5
$ws.Range("L126").Value = 14998.5
$ws.Range("N126").Value = -19938.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4008.3635
$ws.Range("I62").Value = 3542
$ws.Range("K62").Value = 3542
$ws.Range("M62").Value = -2918

$ws.Range("H65").Value = 4008.3635
$ws.Range("I65").Value = 3542
$ws.Range("K65").Value = 17710
$ws.Range("M65").Value = -14590

$ws.Range("H81").Value = 2174.25
$ws.Range("I81").Value = 2174.25
$ws.Range("K81").Value = 4348.5
$ws.Range("M81").Value = -3287.5

$ws.Range("H84").Value = 2174.25
$ws.Range("I84").Value = 2174.25
$ws.Range("K84").Value = 21742.5
$ws.Range("M84").Value = -16438.5

$ws.Range("H122").Value = 1497.5
$ws.Range("I122").Value = 1497.5
$ws.Range("K122").Value = 4492.5
$ws.Range("M122").Value = -2042.5

$ws.Range("H126").Value = 1736.6666
$ws.Range("J126").Value = 1285
$ws.Range("L126").Value = 3855
$ws.Range("N126").Value = -8795

$ws.Range("H132").Value = 1384.5834
$ws.Range("I132").Value = 1384.5834
$ws.Range("K132").Value = 4153.7502
$ws.Range("M132").Value = -1623.7502

$ws.Range("H136").Value = 1377.4166
$ws.Range("I136").Value = 1377.4166
$ws.Range("K136").Value = 4132.2498
$ws.Range("M136").Value = -1582.2498
